$d = $word.ActiveDocument

$d.Content.Find.Execute("David Compernisae", $true, $false, $false, $false, $false, $true, 1, $false, "TRẦN THANH TÂM", 2)
$d.Content.Find.Execute("Califoniare, 198/C, Holiel, Cmalo JST", $true, $false, $false, $false, $false, $true, 1, $false, "Ấp Hoà Tây A, Phú Thuận, Thoại Sơn, An Giang", 2)
$d.Content.Find.Execute("0986545448", $true, $false, $false, $false, $false, $true, 1, $false, "0911938971", 2)
$d.Content.Find.Execute("Accent 1.2 MT 2022", $true, $false, $false, $false, $false, $true, 1, $false, "ACCENT 1.4 AT TIÊU CHUẨN 2021", 2)
$d.Content.Find.Execute("VIN32108412421", $true, $false, $false, $false, $false, $true, 1, $false, "RLUAC41BBPN100453", 2)
$d.Content.Find.Execute("MAY321321", $true, $false, $false, $false, $false, $true, 1, $false, "G4LCPU887669", 2)
